# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
# (and the Quant/Aave + BabyDogeCoin/RenderToken row-order swap) from the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.284.44"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.24%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.871.00"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.29%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7088"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.15%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.59"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.05%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07807"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.93%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3093"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.66%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.00"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.36%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08402"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.39%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.880.92"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.57%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.33%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7099"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.09%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.00"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.40%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.296.64"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.25%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.066"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.00%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008195"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.67%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "239.64"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.69%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.95%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.122.51"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.46%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.752"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.47%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.08%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1591"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.21%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.28"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.10%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.997"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.44%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.43"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.36%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.503"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.19%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.388"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.24%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.295"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.25%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.288"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.72%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05371"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.19%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.938"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.31%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.74%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7471"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.72%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.697"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.51%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01870"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.64%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.229.36"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.87%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.722"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.55%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.549"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.68%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8861"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.31%  "

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.38"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.19%  "

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "108.93"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.33%  "

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.06%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.011.76"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.02%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5192"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.85%  "

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.789"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.62%  "

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000122"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.20%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.403"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.75%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4308"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.33%  "
